$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.718.78"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "1.726.90"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9969"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.96"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9975"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4920"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2616"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06229"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").Value = "1.731.77"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.85"
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06999"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6103"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.504"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.20"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9972"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "26.511.68"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9965"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007205"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.42"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "1.946.54"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.442"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.577"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.107"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.10"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.38"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.389"
$ws.Range("E27").Value = "  -3.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.748"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.38"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.919"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07988"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04498"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.615"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.003"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6261"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9364"
$ws.Range("E37").Value = "  +3.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.009"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.422"
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9973"
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01514"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.581"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.56"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3861"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.907"
$ws.Range("E45").Value = "  +3.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1160"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05384"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.822"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.69"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.236"
$ws.Range("E51").Value = "  +0.05%  "
